# Update reference to latest design, fix BOM
# - Rename the existing "BOM" sheet to "BOM_BOTTOM (default)"
# - Duplicate it as a new sheet "BOM_TOP" (placed right after)
# - On "BOM_TOP", update the J1..J19 header-connector row (row 14) to the
#   newer part (Sullins PPPC151LFBN-RC instead of Adam Tech PH1-15-UA)
# - Restore/update cursor positions on both sheets

$wb = $excel.ActiveWorkbook

$wsBottom = $wb.Worksheets.Item(1)
$wsBottom.Name = "BOM_BOTTOM (default)"

# Duplicate the sheet right after itself to create BOM_TOP
$wsBottom.Copy($null, $wsBottom)
$wsTop = $wb.Worksheets.Item(2)
$wsTop.Name = "BOM_TOP"

# Fix the BOM on BOM_TOP: row 14 (REFDES J1;J2;...;J19) now references the
# updated connector part
$wsTop.Range("D14").Value = "CONN HDR 15POS 0.1 GOLD PCB"
$wsTop.Range("F14").Value = "Sullins"
$wsTop.Range("G14").Value = "PPPC151LFBN-RC"
$wsTop.Range("H14").Value = "S7048-ND"

# Restore the selection/zoom on BOM_TOP
[void]$wsTop.Activate()
[void]$wsTop.Range("C19").Select()
$excel.ActiveWindow.Zoom = 95

# Put the cursor on BOM_BOTTOM (default) and make it the active sheet/tab
[void]$wsBottom.Activate()
[void]$wsBottom.Range("C13").Select()
$excel.ActiveWindow.Zoom = 110
